$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-16 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-17 Friday", 2)

$d.Content.Find.Execute("46×39=", $true, $false, $false, $false, $false, $true, 1, $false, "66×25=", 2)
$d.Content.Find.Execute("21×32=", $true, $false, $false, $false, $false, $true, 1, $false, "98×27=", 2)
$d.Content.Find.Execute("33×57=", $true, $false, $false, $false, $false, $true, 1, $false, "76×27=", 2)
$d.Content.Find.Execute("28×87=", $true, $false, $false, $false, $false, $true, 1, $false, "41×29=", 2)
$d.Content.Find.Execute("32×46=", $true, $false, $false, $false, $false, $true, 1, $false, "70×12=", 2)
$d.Content.Find.Execute("84×34=", $true, $false, $false, $false, $false, $true, 1, $false, "34×54=", 2)
$d.Content.Find.Execute("33×40=", $true, $false, $false, $false, $false, $true, 1, $false, "98×52=", 2)
$d.Content.Find.Execute("47×34=", $true, $false, $false, $false, $false, $true, 1, $false, "75×34=", 2)
$d.Content.Find.Execute("25×23=", $true, $false, $false, $false, $false, $true, 1, $false, "44×71=", 2)
$d.Content.Find.Execute("33×37=", $true, $false, $false, $false, $false, $true, 1, $false, "80×67=", 2)
$d.Content.Find.Execute("56×50=", $true, $false, $false, $false, $false, $true, 1, $false, "29×18=", 2)
$d.Content.Find.Execute("79×34=", $true, $false, $false, $false, $false, $true, 1, $false, "96×62=", 2)
$d.Content.Find.Execute("18×88=", $true, $false, $false, $false, $false, $true, 1, $false, "30×51=", 2)
$d.Content.Find.Execute("45×18=", $true, $false, $false, $false, $false, $true, 1, $false, "75×45=", 2)
$d.Content.Find.Execute("21×60=", $true, $false, $false, $false, $false, $true, 1, $false, "66×95=", 2)
$d.Content.Find.Execute("30×61=", $true, $false, $false, $false, $false, $true, 1, $false, "57×17=", 2)
$d.Content.Find.Execute("43×47=", $true, $false, $false, $false, $false, $true, 1, $false, "16×91=", 2)
$d.Content.Find.Execute("90×69=", $true, $false, $false, $false, $false, $true, 1, $false, "88×81=", 2)
$d.Content.Find.Execute("87×70=", $true, $false, $false, $false, $false, $true, 1, $false, "13×39=", 2)
$d.Content.Find.Execute("19×57=", $true, $false, $false, $false, $false, $true, 1, $false, "68×16=", 2)
$d.Content.Find.Execute("79×64=", $true, $false, $false, $false, $false, $true, 1, $false, "30×64=", 2)
$d.Content.Find.Execute("19×49=", $true, $false, $false, $false, $false, $true, 1, $false, "99×43=", 2)
$d.Content.Find.Execute("12×37=", $true, $false, $false, $false, $false, $true, 1, $false, "66×32=", 2)
$d.Content.Find.Execute("41×34=", $true, $false, $false, $false, $false, $true, 1, $false, "46×70=", 2)
$d.Content.Find.Execute("61×41=", $true, $false, $false, $false, $false, $true, 1, $false, "63×97=", 2)
